$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" (column D) cells hold numeric-looking text
# (e.g. "226.60", "11.06") that must stay text, matching the source data.
# A leading apostrophe forces Excel to keep the literal text instead of
# auto-converting it to a number (which would drop the trailing zero).

$ws.Range('D2').Value = '34.204.61'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').Value = '1.784.03'
$ws.Range('E3').Value = '  +0.17%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '''226.60'
$ws.Range('E5').Value = '  +0.84%  '

$ws.Range('E6').Value = '  +0.34%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = '''31.94'
$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('E9').Value = '  +0.87%  '

$ws.Range('E10').Value = '  +2.19%  '

$ws.Range('D11').Value = '''0.0947'
$ws.Range('E11').Value = '  +1.16%  '

$ws.Range('D12').Value = '2.042.04'
$ws.Range('E12').Value = '  +0.25%  '

$ws.Range('D13').Value = '''11.06'
$ws.Range('E13').Value = '  -1.83%  '

$ws.Range('D14').Value = '1.795.07'
$ws.Range('E14').Value = '  +0.79%  '

$ws.Range('E15').Value = '  +2.30%  '

$ws.Range('D16').Value = '34.181.51'
$ws.Range('E16').Value = '  +0.85%  '

$ws.Range('D17').Value = '''4.19'
$ws.Range('E17').Value = '  +1.19%  '

$ws.Range('D18').Value = '''67.99'
$ws.Range('E18').Value = '  +1.92%  '

$ws.Range('D19').Value = '0.0₃0808'
$ws.Range('E19').Value = '  +4.47%  '

$ws.Range('D20').Value = '''247.17'
$ws.Range('E20').Value = '  +3.54%  '

$ws.Range('E21').Value = '  +4.00%  '

$ws.Range('E22').Value = '  +0.11%  '

$ws.Range('E23').Value = '  +2.11%  '

$ws.Range('E24').Value = '  -1.23%  '

$ws.Range('D25').Value = '''162.71'
$ws.Range('E25').Value = '  +1.33%  '

$ws.Range('E26').Value = '  +2.29%  '

$ws.Range('D27').Value = '''16.31'

$ws.Range('E28').Value = '  +1.59%  '

$ws.Range('E30').Value = '  +0.60%  '

$ws.Range('E31').Value = '  +1.98%  '

$ws.Range('D32').Value = '''3.75'
$ws.Range('E32').Value = '  +4.23%  '

$ws.Range('D33').Value = '''3.76'
$ws.Range('E33').Value = '  +6.67%  '

$ws.Range('E34').Value = '  -1.46%  '

$ws.Range('D35').Value = '1.444.30'
$ws.Range('E35').Value = '  +3.91%  '

$ws.Range('E36').Value = '  +2.78%  '

$ws.Range('D37').Value = '''2.40'
$ws.Range('E37').Value = '  +6.49%  '

$ws.Range('E38').Value = '  +3.64%  '

$ws.Range('E39').Value = '  +0.62%  '

$ws.Range('D40').Value = '''80.36'
$ws.Range('E40').Value = '  +2.18%  '

$ws.Range('D41').Value = '''2.36'
$ws.Range('E41').Value = '  -1.03%  '

$ws.Range('D42').Value = '''0.926'
$ws.Range('E42').Value = '  +1.57%  '

$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.67'
$ws.Range('E43').Value = '  +0.80%  '

$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '''13.54'
$ws.Range('E44').Value = '  -0.56%  '

$ws.Range('E45').Value = '  +3.77%  '

$ws.Range('E46').Value = '  +0.76%  '

$ws.Range('E47').Value = '  -0.16%  '

$ws.Range('D48').Value = '0.0₆0134'
$ws.Range('E48').Value = '  -3.87%  '

$ws.Range('D49').Value = '1.943.35'
$ws.Range('E49').Value = '  +0.26%  '

$ws.Range('D50').Value = '''104.67'
$ws.Range('E50').Value = '  -2.03%  '

$ws.Range('E51').Value = '  +0.09%  '
